$wb = $excel.ActiveWorkbook

# Report regenerated for handoff: status flips from "In Translation" to
# "Ready for handoff" and the associated timestamps advance ~40s. The wider
# status text means the "Status" columns need to grow to fit it.
$newWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-27 02:38:36"
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-27 02:38:32"
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-27 02:38:36"
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
